# "fix: remove duplicated cover" - the deck has a single, duplicated
# cover slide (empty Title + Text Placeholder) left over; drop it so the
# presentation ends up with zero slides, matching the canonical edit.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Delete()
